$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.873717665672302
$ws.Range("B1").Value = 4.774370193481445
$ws.Range("C1").Value = 3.642581939697266
$ws.Range("D1").Value = 1.212711691856384
$ws.Range("E1").Value = 0.7983418703079224
